# Auto-generated script applying numeric updates to Gilgamesh Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 11906533
$ws.Range("J17").Value = 11906533
$ws.Range("L17").Value = 35719599
$ws.Range("N17").Value = -35719935

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 156.25
$ws.Range("I33").Value = 162.36363
$ws.Range("J33").Value = 89
$ws.Range("K33").Value = 162.36363
$ws.Range("L33").Value = 89
$ws.Range("M33").Value = 66.63637
$ws.Range("N33").Value = -547

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 500010000
$ws.Range("I64").Value = 20000
$ws.Range("K64").Value = 20000
$ws.Range("M64").Value = -19752

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 500010000
$ws.Range("I67").Value = 20000
$ws.Range("K67").Value = 20000
$ws.Range("M67").Value = -19142

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2974.2327
$ws.Range("I98").Value = 2974.2327
$ws.Range("K98").Value = 2974.2327
$ws.Range("M98").Value = -1476.2327

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2974.2327
$ws.Range("I122").Value = 2974.2327
$ws.Range("K122").Value = 8922.6981
$ws.Range("M122").Value = -6472.6981

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 5375.346
$ws.Range("J137").Value = 10354.728
$ws.Range("L137").Value = 31064.184
$ws.Range("N137").Value = -36164.18399999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 4212.875
$ws.Range("J141").Value = 4726.25
$ws.Range("L141").Value = 14178.75
$ws.Range("N141").Value = -24538.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 559.8
$ws.Range("I4").Value = 466.33334
$ws.Range("K4").Value = 466.33334
$ws.Range("M4").Value = -350.33334

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 29172.842
$ws.Range("I45").Value = 38871.54
$ws.Range("J45").Value = 8159
$ws.Range("K45").Value = 38871.54
$ws.Range("L45").Value = 8159
$ws.Range("M45").Value = -38494.54
$ws.Range("N45").Value = -8913

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H50").Value = 415.75
$ws.Range("I50").Value = 0
$ws.Range("J50").Value = 415.75
$ws.Range("K50").Value = 0
$ws.Range("L50").Value = 415.75
$ws.Range("M50").ClearContents()
$ws.Range("N50").Value = -1843.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4549.875
$ws.Range("I61").Value = 1797.875
$ws.Range("K61").Value = 1797.875
$ws.Range("M61").Value = -1585.875

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 1616.619
$ws.Range("I97").Value = 1572.5
$ws.Range("K97").Value = 1572.5
$ws.Range("M97").Value = -1076.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 2999.9614
$ws.Range("I122").Value = 2614.2766
$ws.Range("J122").Value = 6625.4
$ws.Range("K122").Value = 7842.8298
$ws.Range("L122").Value = 19876.2
$ws.Range("M122").Value = -5392.8298
$ws.Range("N122").Value = -24776.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H133").Value = 66000
$ws.Range("J133").Value = 66000
$ws.Range("L133").Value = 66000
$ws.Range("N133").Value = -71060

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4549.875
$ws.Range("I136").Value = 1797.875
$ws.Range("K136").Value = 5393.625
$ws.Range("M136").Value = -2843.625

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4224.5
$ws.Range("I86").Value = 4121.4443
$ws.Range("J86").Value = 4357
$ws.Range("K86").Value = 4121.4443
$ws.Range("L86").Value = 4357
$ws.Range("M86").Value = -2998.4443
$ws.Range("N86").Value = -6603

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 4224.5
$ws.Range("I89").Value = 4121.4443
$ws.Range("J89").Value = 4357
$ws.Range("K89").Value = 20607.2215
$ws.Range("L89").Value = 21785
$ws.Range("M89").Value = -14991.2215
$ws.Range("N89").Value = -33017

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2382.5625
$ws.Range("I16").Value = 2370.3635
$ws.Range("J16").Value = 2409.4
$ws.Range("K16").Value = 2370.3635
$ws.Range("L16").Value = 2409.4
$ws.Range("M16").Value = -2083.3635
$ws.Range("N16").Value = -2983.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4836.967
$ws.Range("J31").Value = 4903.5
$ws.Range("L31").Value = 4903.5
$ws.Range("N31").Value = -5493.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4836.967
$ws.Range("J34").Value = 4903.5
$ws.Range("L34").Value = 4903.5
$ws.Range("N34").Value = -5307.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H68").Value = 88372
$ws.Range("J68").Value = 88372
$ws.Range("L68").Value = 88372
$ws.Range("N68").Value = -89870

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H71").Value = 88372
$ws.Range("J71").Value = 88372
$ws.Range("L71").Value = 265116
$ws.Range("N71").Value = -272604

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 2382.5625
$ws.Range("I113").Value = 2370.3635
$ws.Range("J113").Value = 2409.4
$ws.Range("K113").Value = 2370.3635
$ws.Range("L113").Value = 2409.4
$ws.Range("M113").Value = -200.3634999999999
$ws.Range("N113").Value = -6749.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 3470.0557
$ws.Range("I122").Value = 2823.4666
$ws.Range("J122").Value = 6703
$ws.Range("K122").Value = 8470.399800000001
$ws.Range("L122").Value = 20109
$ws.Range("M122").Value = -6020.399800000001
$ws.Range("N122").Value = -25009

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 29416388
$ws.Range("I132").Value = 45458330
$ws.Range("K132").Value = 136374990
$ws.Range("M132").Value = -136372460

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 5874.222
$ws.Range("I134").Value = 5808
$ws.Range("K134").Value = 17424
$ws.Range("M134").Value = -14889

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 6921.0557
$ws.Range("I56").Value = 6921.0557
$ws.Range("K56").Value = 6921.0557
$ws.Range("M56").Value = -6391.0557

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H96").Value = 4718.75
$ws.Range("J96").Value = 4718.75
$ws.Range("L96").Value = 14156.25
$ws.Range("N96").Value = -18274.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H109").Value = 2027
$ws.Range("I109").Value = 796.2857
$ws.Range("J109").Value = 3750
$ws.Range("K109").Value = 2388.8571
$ws.Range("L109").Value = 11250
$ws.Range("M109").Value = -1348.8571
$ws.Range("N109").Value = -13330

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6308.6665
$ws.Range("J113").Value = 6920.316
$ws.Range("L113").Value = 20760.948
$ws.Range("N113").Value = -25100.948

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H115").Value = 220107.08
$ws.Range("I115").Value = 2000
$ws.Range("J115").Value = 236884.53
$ws.Range("K115").Value = 6000
$ws.Range("L115").Value = 710653.59
$ws.Range("M115").Value = -4825
$ws.Range("N115").Value = -713003.59

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 2267.5
$ws.Range("I131").Value = 2223.5557
$ws.Range("K131").Value = 6670.6671
$ws.Range("M131").Value = -1630.6671

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 773.8182
$ws.Range("I2").Value = 946.125
$ws.Range("K2").Value = 946.125
$ws.Range("M2").Value = -833.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2596.2778
$ws.Range("I132").Value = 2211.3635
$ws.Range("J132").Value = 3201.1428
$ws.Range("K132").Value = 6634.0905
$ws.Range("L132").Value = 9603.428400000001
$ws.Range("M132").Value = -4104.0905
$ws.Range("N132").Value = -14663.4284

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H135").Value = 68298.87
$ws.Range("J135").Value = 68298.87
$ws.Range("L135").Value = 68298.87
$ws.Range("N135").Value = -78438.87

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H137").Value = 77777
$ws.Range("J137").Value = 77777
$ws.Range("L137").Value = 77777
$ws.Range("N137").Value = -87977

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H140").Value = 82554.78
$ws.Range("J140").Value = 82554.78
$ws.Range("L140").Value = 82554.78
$ws.Range("N140").Value = -92914.78

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H141").Value = 68997.336
$ws.Range("J141").Value = 68997.336
$ws.Range("L141").Value = 68997.336
$ws.Range("N141").Value = -79357.336

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 411.9
$ws.Range("I22").Value = 448.75
$ws.Range("J22").Value = 264.5
$ws.Range("K22").Value = 448.75
$ws.Range("L22").Value = 264.5
$ws.Range("M22").Value = -153.75
$ws.Range("N22").Value = -854.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 411.9
$ws.Range("I27").Value = 448.75
$ws.Range("J27").Value = 264.5
$ws.Range("K27").Value = 448.75
$ws.Range("L27").Value = 264.5
$ws.Range("M27").Value = -341.75
$ws.Range("N27").Value = -478.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 27117.25
$ws.Range("I40").Value = 29982.896
$ws.Range("K40").Value = 29982.896
$ws.Range("M40").Value = -29846.896

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2117.6
$ws.Range("I122").Value = 2071.2104
$ws.Range("J122").Value = 2999
$ws.Range("K122").Value = 6213.6312
$ws.Range("L122").Value = 8997
$ws.Range("M122").Value = -3763.6312
$ws.Range("N122").Value = -13897

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H140").Value = 78194.75
$ws.Range("J140").Value = 78194.75
$ws.Range("L140").Value = 78194.75
$ws.Range("N140").Value = -88554.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H56").Value = 74998
$ws.Range("J56").Value = 74998
$ws.Range("L56").Value = 74998
$ws.Range("N56").Value = -76426

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 10420115
$ws.Range("I132").Value = 13336868
$ws.Range("J132").Value = 3142.5715
$ws.Range("K132").Value = 40010604
$ws.Range("L132").Value = 9427.7145
$ws.Range("M132").Value = -40008074
$ws.Range("N132").Value = -14487.7145
